$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# NumberFormat "@" + Style reset keeps numeric-looking price strings as text
# (matching the source feed's inlineStr cells) without leaving a residual style on the cell.

$d = $ws.Range("D2")
$d.NumberFormat = "@"
$d.Value = '28.918.89'
$d.Style = "Normal"
$ws.Range("E2").Value = '  +1.27%  '

$d = $ws.Range("D3")
$d.NumberFormat = "@"
$d.Value = '1.878.99'
$d.Style = "Normal"
$ws.Range("E3").Value = '  -0.68%  '

$d = $ws.Range("D4")
$d.NumberFormat = "@"
$d.Value = '1.001'
$d.Style = "Normal"
$ws.Range("E4").Value = '  -0.76%  '

$d = $ws.Range("D5")
$d.NumberFormat = "@"
$d.Value = '325.05'
$d.Style = "Normal"
$ws.Range("E5").Value = '  -0.54%  '

$ws.Range("E6").Value = '  -0.68%  '

$d = $ws.Range("D7")
$d.NumberFormat = "@"
$d.Value = '0.4595'
$d.Style = "Normal"
$ws.Range("E7").Value = '  -0.03%  '

$d = $ws.Range("D8")
$d.NumberFormat = "@"
$d.Value = '0.3877'
$d.Style = "Normal"
$ws.Range("E8").Value = '  +0.19%  '

$d = $ws.Range("D9")
$d.NumberFormat = "@"
$d.Value = '0.07870'
$d.Style = "Normal"
$ws.Range("E9").Value = '  -0.06%  '

$d = $ws.Range("D10")
$d.NumberFormat = "@"
$d.Value = '0.9857'
$d.Style = "Normal"
$ws.Range("E10").Value = '  -1.71%  '

$d = $ws.Range("D11")
$d.NumberFormat = "@"
$d.Value = '21.79'
$d.Style = "Normal"
$ws.Range("E11").Value = '  +0.70%  '

$d = $ws.Range("D12")
$d.NumberFormat = "@"
$d.Value = '1.875.49'
$d.Style = "Normal"
$ws.Range("E12").Value = '  -1.77%  '

$d = $ws.Range("D13")
$d.NumberFormat = "@"
$d.Value = '6.982'
$d.Style = "Normal"
$ws.Range("E13").Value = '  -1.53%  '

$d = $ws.Range("D14")
$d.NumberFormat = "@"
$d.Value = '5.648'
$d.Style = "Normal"
$ws.Range("E14").Value = '  -1.19%  '

$d = $ws.Range("D15")
$d.NumberFormat = "@"
$d.Value = '0.06963'
$d.Style = "Normal"
$ws.Range("E15").Value = '  -0.09%  '

$d = $ws.Range("D16")
$d.NumberFormat = "@"
$d.Value = '88.04'
$d.Style = "Normal"
$ws.Range("E16").Value = '  +0.59%  '

$d = $ws.Range("D17")
$d.NumberFormat = "@"
$d.Value = '1.002'
$d.Style = "Normal"
$ws.Range("E17").Value = '  -0.70%  '

$d = $ws.Range("D18")
$d.NumberFormat = "@"
$d.Value = '0.000009971'
$d.Style = "Normal"
$ws.Range("E18").Value = '  -0.88%  '

$d = $ws.Range("D19")
$d.NumberFormat = "@"
$d.Value = '16.97'
$d.Style = "Normal"
$ws.Range("E19").Value = '  -1.40%  '

$d = $ws.Range("D20")
$d.NumberFormat = "@"
$d.Value = '1.002'
$d.Style = "Normal"
$ws.Range("E20").Value = '  -0.48%  '

$d = $ws.Range("D21")
$d.NumberFormat = "@"
$d.Value = '28.923.84'
$d.Style = "Normal"
$ws.Range("E21").Value = '  +1.15%  '

$d = $ws.Range("D22")
$d.NumberFormat = "@"
$d.Value = '5.243'
$d.Style = "Normal"
$ws.Range("E22").Value = '  -1.58%  '

$d = $ws.Range("D23")
$d.NumberFormat = "@"
$d.Value = '10.96'
$d.Style = "Normal"
$ws.Range("E23").Value = '  -0.55%  '

$d = $ws.Range("D24")
$d.NumberFormat = "@"
$d.Value = '2.102'
$d.Style = "Normal"
$ws.Range("E24").Value = '  +1.94%  '

$d = $ws.Range("D25")
$d.NumberFormat = "@"
$d.Value = '156.18'
$d.Style = "Normal"
$ws.Range("E25").Value = '  +0.78%  '

$d = $ws.Range("D26")
$d.NumberFormat = "@"
$d.Value = '19.30'
$d.Style = "Normal"

$d = $ws.Range("D27")
$d.NumberFormat = "@"
$d.Value = '6.017'
$d.Style = "Normal"
$ws.Range("E27").Value = '  +2.80%  '

$d = $ws.Range("D28")
$d.NumberFormat = "@"
$d.Value = '1.930'
$d.Style = "Normal"
$ws.Range("E28").Value = '  -1.71%  '

$d = $ws.Range("D29")
$d.NumberFormat = "@"
$d.Value = '117.32'
$d.Style = "Normal"
$ws.Range("E29").Value = '  -0.89%  '

$d = $ws.Range("D30")
$d.NumberFormat = "@"
$d.Value = '0.09321'
$d.Style = "Normal"
$ws.Range("E30").Value = '  -0.13%  '

$d = $ws.Range("D31")
$d.NumberFormat = "@"
$d.Value = '0.9023'
$d.Style = "Normal"
$ws.Range("E31").Value = '  -2.55%  '

$d = $ws.Range("D32")
$d.NumberFormat = "@"
$d.Value = '5.254'
$d.Style = "Normal"
$ws.Range("E32").Value = '  -0.85%  '

$ws.Range("E33").Value = '  -1.44%  '

$d = $ws.Range("D34")
$d.NumberFormat = "@"
$d.Value = '3.254'
$d.Style = "Normal"
$ws.Range("E34").Value = '  -0.47%  '

$d = $ws.Range("D35")
$d.NumberFormat = "@"
$d.Value = '1.186'
$d.Style = "Normal"
$ws.Range("E35").Value = '  +2.66%  '

$d = $ws.Range("D36")
$d.NumberFormat = "@"
$d.Value = '0.05754'
$d.Style = "Normal"
$ws.Range("E36").Value = '  -0.29%  '

$d = $ws.Range("D37")
$d.NumberFormat = "@"
$d.Value = '0.02070'
$d.Style = "Normal"
$ws.Range("E37").Value = '  -0.14%  '

$d = $ws.Range("D38")
$d.NumberFormat = "@"
$d.Value = '1.001'
$d.Style = "Normal"
$ws.Range("E38").Value = '  -0.55%  '

$ws.Range("E39").Value = '  -1.90%  '

$d = $ws.Range("D40")
$d.NumberFormat = "@"
$d.Value = '0.5648'
$d.Style = "Normal"
$ws.Range("E40").Value = '  -0.43%  '

$d = $ws.Range("D41")
$d.NumberFormat = "@"
$d.Value = '0.1764'
$d.Style = "Normal"
$ws.Range("E41").Value = '  -1.51%  '

$d = $ws.Range("D42")
$d.NumberFormat = "@"
$d.Value = '9.659'
$d.Style = "Normal"

$d = $ws.Range("D43")
$d.NumberFormat = "@"
$d.Value = '2.258'
$d.Style = "Normal"
$ws.Range("E43").Value = '  +2.67%  '

$ws.Range("E44").Value = '  +0.47%  '

$d = $ws.Range("D45")
$d.NumberFormat = "@"
$d.Value = '0.5345'
$d.Style = "Normal"
$ws.Range("E45").Value = '  -0.30%  '

$d = $ws.Range("D46")
$d.NumberFormat = "@"
$d.Value = '0.07044'
$d.Style = "Normal"
$ws.Range("E46").Value = '  -1.65%  '

$d = $ws.Range("D47")
$d.NumberFormat = "@"
$d.Value = '1.847'
$d.Style = "Normal"
$ws.Range("E47").Value = '  +0.24%  '

$d = $ws.Range("D48")
$d.NumberFormat = "@"
$d.Value = '113.01'
$d.Style = "Normal"
$ws.Range("E48").Value = '  +0.15%  '

$ws.Range("E49").Value = '  +1.58%  '

$ws.Range("E50").Value = '  -5.18%  '

$d = $ws.Range("D51")
$d.NumberFormat = "@"
$d.Value = '70.58'
$d.Style = "Normal"
$ws.Range("E51").Value = '  -0.56%  '

